$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 86, shifting the existing rows 86:194 down to 87:195.
$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with the new daily price record.
$ws.Range("A86").Value = 3
$ws.Range("B86").Value = "Femacal de La Calera"
$ws.Range("C86").Value = "Coquimbo"
$ws.Range("D86").Value = 44467
$ws.Range("E86").Value = 5
$ws.Range("F86").Value = 100114013
$ws.Range("G86").Value = "Zanahoria"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 510
$ws.Range("K86").Value = 6500
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = 6725
$ws.Range("N86").Value = "$/saco 20 kilos"
$ws.Range("O86").Value = "Chillán"
$ws.Range("P86").Value = 336
$ws.Range("Q86").Value = 20
$ws.Range("R86").Value = "Hortaliza"
